$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2021
$ws.Range("I62").Value = 1302.5
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 1302.5
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -678.5
$ws.Range("N62").Value = -3748

$ws.Range("H65").Value = 2021
$ws.Range("I65").Value = 1302.5
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 6512.5
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -3392.5
$ws.Range("N65").Value = -18740

$ws.Range("H68").Value = 28495
$ws.Range("J68").Value = 28495
$ws.Range("L68").Value = 28495
$ws.Range("N68").Value = -29993

$ws.Range("H71").Value = 28495
$ws.Range("J71").Value = 28495
$ws.Range("L71").Value = 85485
$ws.Range("N71").Value = -92973

$ws.Range("H132").Value = 37720560
$ws.Range("I132").Value = 46237464
$ws.Range("J132").Value = 2836.5715
$ws.Range("K132").Value = 138712392
$ws.Range("L132").Value = 8509.7145
$ws.Range("M132").Value = -138709862
$ws.Range("N132").Value = -13569.7145

$ws.Range("H137").Value = 794754.3
$ws.Range("I137").Value = 1191228.5
$ws.Range("K137").Value = 3573685.5
$ws.Range("M137").Value = -3571135.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1005.65
$ws.Range("I2").Value = 1005.65
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1005.65
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -892.65
$ws.Range("N2").Value = $null

$ws.Range("H45").Value = 1150.4375
$ws.Range("I45").Value = 902.3333
$ws.Range("J45").Value = 1207.6923
$ws.Range("K45").Value = 902.3333
$ws.Range("L45").Value = 1207.6923
$ws.Range("M45").Value = -525.3333
$ws.Range("N45").Value = -1961.6923

$ws.Range("H50").Value = 867
$ws.Range("I50").Value = 1062.25
$ws.Range("J50").Value = 476.5
$ws.Range("K50").Value = 1062.25
$ws.Range("L50").Value = 476.5
$ws.Range("M50").Value = -348.25
$ws.Range("N50").Value = -1904.5

$ws.Range("H95").Value = 19998.2
$ws.Range("J95").Value = 19998.2
$ws.Range("L95").Value = 19998.2
$ws.Range("N95").Value = -25490.2

$ws.Range("H116").Value = 1005.65
$ws.Range("I116").Value = 1005.65
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1005.65
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1288.35
$ws.Range("N116").Value = $null

$ws.Range("H132").Value = 2660872.2
$ws.Range("I132").Value = 3473050
$ws.Range("J132").Value = 2835.2727
$ws.Range("K132").Value = 10419150
$ws.Range("L132").Value = 8505.8181
$ws.Range("M132").Value = -10416620
$ws.Range("N132").Value = -13565.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1005.65
$ws.Range("I3").Value = 1005.65
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1005.65
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -891.65
$ws.Range("N3").Value = $null

$ws.Range("H86").Value = 1335.1111
$ws.Range("I86").Value = 1403.2
$ws.Range("K86").Value = 1403.2
$ws.Range("M86").Value = -280.2

$ws.Range("H89").Value = 1335.1111
$ws.Range("I89").Value = 1403.2
$ws.Range("K89").Value = 7016
$ws.Range("M89").Value = -1400

$ws.Range("H99").Value = 818.5
$ws.Range("I99").Value = 818.75
$ws.Range("J99").Value = 818.1667
$ws.Range("K99").Value = 818.75
$ws.Range("L99").Value = 818.1667
$ws.Range("M99").Value = 679.25
$ws.Range("N99").Value = -3814.1667

$ws.Range("H134").Value = 10117585
$ws.Range("I134").Value = 11512790
$ws.Range("K134").Value = 34538370
$ws.Range("M134").Value = -34535835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16565.852
$ws.Range("I31").Value = 1197
$ws.Range("J31").Value = 28860.934
$ws.Range("K31").Value = 1197
$ws.Range("L31").Value = 28860.934
$ws.Range("M31").Value = -902
$ws.Range("N31").Value = -29450.934

$ws.Range("H34").Value = 16565.852
$ws.Range("I34").Value = 1197
$ws.Range("J34").Value = 28860.934
$ws.Range("K34").Value = 1197
$ws.Range("L34").Value = 28860.934
$ws.Range("M34").Value = -995
$ws.Range("N34").Value = -29264.934

$ws.Range("H122").Value = 5756410.5
$ws.Range("I122").Value = 12988568
$ws.Range("J122").Value = 74001
$ws.Range("K122").Value = 38965704
$ws.Range("L122").Value = 222003
$ws.Range("M122").Value = -38963254
$ws.Range("N122").Value = -226903

$ws.Range("H132").Value = 8550844
$ws.Range("I132").Value = 11495052
$ws.Range("J132").Value = 12640.2
$ws.Range("K132").Value = 34485156
$ws.Range("L132").Value = 37920.60000000001
$ws.Range("M132").Value = -34482626
$ws.Range("N132").Value = -42980.60000000001

$ws.Range("H134").Value = 11162091
$ws.Range("I134").Value = 13890128
$ws.Range("J134").Value = 6251626
$ws.Range("K134").Value = 41670384
$ws.Range("L134").Value = 18754878
$ws.Range("M134").Value = -41667849
$ws.Range("N134").Value = -18759948

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 49742.75
$ws.Range("J101").Value = 49742.75
$ws.Range("L101").Value = 49742.75
$ws.Range("N101").Value = -56232.75

$ws.Range("H107").Value = 213
$ws.Range("I107").Value = 193.73334
$ws.Range("J107").Value = 254.28572
$ws.Range("K107").Value = 193.73334
$ws.Range("L107").Value = 254.28572
$ws.Range("M107").Value = 1726.26666
$ws.Range("N107").Value = -4094.28572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1667.84
$ws.Range("I7").Value = 1456.125
$ws.Range("J7").Value = 2044.2222
$ws.Range("K7").Value = 1456.125
$ws.Range("L7").Value = 2044.2222
$ws.Range("M7").Value = -1344.125
$ws.Range("N7").Value = -2268.2222

$ws.Range("H122").Value = 78127300
$ws.Range("I122").Value = 125001640
$ws.Range("J122").Value = 31252950
$ws.Range("K122").Value = 375004920
$ws.Range("L122").Value = 93758850
$ws.Range("M122").Value = -375002470
$ws.Range("N122").Value = -93763750

$ws.Range("H126").Value = 1667.84
$ws.Range("I126").Value = 1456.125
$ws.Range("J126").Value = 2044.2222
$ws.Range("K126").Value = 4368.375
$ws.Range("L126").Value = 6132.6666
$ws.Range("M126").Value = -1898.375
$ws.Range("N126").Value = -11072.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15383.429
$ws.Range("I122").Value = 17667.334
$ws.Range("J122").Value = 1680
$ws.Range("K122").Value = 53002.00199999999
$ws.Range("L122").Value = 5040
$ws.Range("M122").Value = -50552.00199999999
$ws.Range("N122").Value = -9940
